$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "2022-Q3" sheet, inserted right before the current "2022-Q2"
#    sheet (which sits at position 2, right after "总计"). Copying the
#    "2022-Q2" sheet gives us an exact structural/style clone (sheetPr,
#    column styles, header formatting, page margins, etc.) that we then
#    overwrite with the 2022-Q3 figures.
# ---------------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item(2)
$sheetQ2.Copy($sheetQ2) | Out-Null
$sheetQ3 = $wb.Worksheets.Item(2)
$sheetQ3.Name = "2022-Q3"

# Row 2 - 008715 / 景顺长城价值驱动一年持有期灵活配置混合
$sheetQ3.Range("B2:G2").NumberFormat = "@"
$sheetQ3.Range("B2").Value = "008715"
$sheetQ3.Range("C2").Value = "景顺长城价值驱动一年持有期灵活配置混合"
$sheetQ3.Range("D2").Value = "8.44"
$sheetQ3.Range("E2").Value = "90.91"
$sheetQ3.Range("F2").Value = "1.43"
$sheetQ3.Range("G2").Value = "0.1207"
$sheetQ3.Range("H2").Value = 10

# Row 3 - 009098 / 景顺长城价值领航两年持有期混合
$sheetQ3.Range("B3:G3").NumberFormat = "@"
$sheetQ3.Range("B3").Value = "009098"
$sheetQ3.Range("C3").Value = "景顺长城价值领航两年持有期混合"
$sheetQ3.Range("D3").Value = "7.16"
$sheetQ3.Range("E3").Value = "90.92"
$sheetQ3.Range("F3").Value = "1.46"
$sheetQ3.Range("G3").Value = "0.1045"
$sheetQ3.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new data row right under the
#    header for the 2022-Q3 totals, pushing every later quarter down by one
#    row (their own figures stay untouched, only their row position moves).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# Re-apply the same formatting the other data rows use for column A, since a
# freshly inserted row starts out blank/unstyled.
$summary.Range("A3:D3").Copy() | Out-Null
$summary.Range("A2:D2").PasteSpecial(-4122) | Out-Null

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.23

# Renumber the 0-based index column for the rows that just shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
